$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily report")

# Fill in the three notes cells for row 9 (task 8) that were previously blank.
# Order matters for shared-string table assignment: E9, then F9, then D9.
$ws.Range("E9").Value = "Experimented with reflection mix levels and volume adjustments to address EDT and RT60 issues. Tested normalized sine sweep with lower amplitude. Compared results with Mona's recordings. Implemented volume attenuation in Unity (-10dB to -25dB). Observed persistent high RT60 due to elevated noise floor levels. Identified long reverb tail in RIR time domain graph. Concluded current approaches ineffective in resolving noise floor problem."
$ws.Range("F9").Value = "Created test scenes to isolate Steam Audio issues. Tested empty scene with floor, confirming proper deconvolution but persistent noise floor. Experimented with floorless scene, revealing unexpected small reverb and noise floor. Observed second peak in RIR not reflected in octave band graph. Planned to replicate findings using Mona's code and test on ground truth and S3A project meshes. Initiated research into Unity RIR measurement practices."
$ws.Range("D9").Value = " Tested reflection mix level changes. Prepared questions for Dr. Atiyeh on RIR analysis. Identified and fixed deconvolution method issue. Began rewriting KT Jupyter notebook for compactness and correctness. Experimented with y_fit parameter, noting persistent EDT issues. Planned further testing and refinement for tomorrow."

# Grow row 9 to fit the new, taller wrapped text (116 -> 174).
$ws.Rows.Item(9).RowHeight = 174

# Move the active selection to B10, matching the post-edit cursor position.
$ws.Range("B10").Select()
